# Apply updated coin price/volume data to worksheet cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe forces Excel to store a numeric-looking string as text,
# matching the inlineStr cell type used in the source workbook.

$ws.Range("D2").Value = '67.300.71'
$ws.Range("E2").Value = '  -1.72%  '

$ws.Range("D3").Value = '3.507.59'
$ws.Range("E3").Value = '  -3.73%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = '''200.37'
$ws.Range("E5").Value = '  +1.76%  '

$ws.Range("D6").Value = '''553.11'
$ws.Range("E6").Value = '  -4.89%  '

$ws.Range("D7").Value = '3.501.61'
$ws.Range("E7").Value = '  -3.68%  '

$ws.Range("E8").Value = '  -2.23%  '

$ws.Range("E9").Value = '  -0.05%  '

$ws.Range("E10").Value = '  -3.75%  '

$ws.Range("D11").Value = '''62.17'
$ws.Range("E11").Value = '  +10.04%  '

$ws.Range("E12").Value = '  -7.21%  '

$ws.Range("E13").Value = '  -7.58%  '

$ws.Range("D14").Value = '''9.84'
$ws.Range("E14").Value = '  -2.64%  '

$ws.Range("D15").Value = '4.064.80'
$ws.Range("E15").Value = '  -3.84%  '

$ws.Range("D16").Value = '3.509.73'
$ws.Range("E16").Value = '  -3.71%  '

$ws.Range("E17").Value = '  -1.94%  '

$ws.Range("D18").Value = '''18.45'
$ws.Range("E18").Value = '  -1.31%  '

$ws.Range("D19").Value = '67.107.47'
$ws.Range("E19").Value = '  -1.93%  '

$ws.Range("D20").Value = '''11.83'
$ws.Range("E20").Value = '  -6.17%  '

$ws.Range("E21").Value = '  -5.48%  '

$ws.Range("D22").Value = '''392.16'
$ws.Range("E22").Value = '  -2.77%  '

$ws.Range("D23").Value = '''12.40'
$ws.Range("E23").Value = '  -4.60%  '

$ws.Range("E24").Value = '  -6.14%  '

$ws.Range("D25").Value = '''83.03'
$ws.Range("E25").Value = '  -3.69%  '

$ws.Range("D26").Value = '''3.95'
$ws.Range("E26").Value = '  +2.22%  '

$ws.Range("D27").Value = '''12.26'
$ws.Range("E27").Value = '  -3.06%  '

$ws.Range("E28").Value = '  -5.07%  '

$ws.Range("D29").Value = '''8.87'
$ws.Range("E29").Value = '  -3.69%  '

$ws.Range("D30").Value = '''31.02'
$ws.Range("E30").Value = '  -2.57%  '

$ws.Range("D31").Value = '''690.24'
$ws.Range("E31").Value = '  -1.69%  '

$ws.Range("D32").Value = '''7.12'
$ws.Range("E32").Value = '  -12.14%  '

$ws.Range("D33").Value = '''11.72'
$ws.Range("E33").Value = '  -4.51%  '

$ws.Range("D34").Value = '''63.78'
$ws.Range("E34").Value = '  -1.70%  '

$ws.Range("E35").Value = '  -6.06%  '

$ws.Range("D36").Value = '''38.77'
$ws.Range("E36").Value = '  -9.52%  '

$ws.Range("E37").Value = '  +0.08%  '

$ws.Range("D38").Value = '''0.398'
$ws.Range("E38").Value = '  -6.49%  '

$ws.Range("E39").Value = '  -5.22%  '

$ws.Range("B40").Value = 'ThetaToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D40").Value = '''3.01'
$ws.Range("E40").Value = '  -4.38%  '

$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '3.065.69'
$ws.Range("E41").Value = '  -5.00%  '

$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '''0.998'
$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("D43").Value = '''2.59'
$ws.Range("E43").Value = '  -9.98%  '

$ws.Range("D44").Value = '0.0₃0678'
$ws.Range("E44").Value = '  -14.13%  '

$ws.Range("D45").Value = '''2.78'
$ws.Range("E45").Value = '  +5.99%  '

$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0403'
$ws.Range("E46").Value = '  -4.66%  '

$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = '''2.68'
$ws.Range("E47").Value = '  -10.70%  '

$ws.Range("D48").Value = '''0.127'
$ws.Range("E48").Value = '  -3.87%  '

$ws.Range("D49").Value = '''138.42'
$ws.Range("E49").Value = '  -2.82%  '

$ws.Range("E50").Value = '  -7.69%  '

$ws.Range("D51").Value = '''2.87'
$ws.Range("E51").Value = '  -7.94%  '
